# Task7-1.xlsx: "Modified Result Task TuTM, ThiVT"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet view: selection moved to K25 ---
# (topLeftCell/zoomScale* aren't round-tripped by this host; selection is)
$ws.Range("K25").Select()

# --- Rows 5 & 6: mark I/J as "needs attention" (red font) ---
$ws.Range("I5:J5").Font.ColorIndex = 3
$ws.Range("I6:J6").Font.ColorIndex = 3

# --- Row 12 (TuTM's task "Study Zend_Auth and make a  LAB") ---
$ws.Range("I12:J12").Font.ColorIndex = 3
$ws.Range("N12").Value = "LinhTA"
$ws.Range("N12").Style = "Normal"
$ws.Range("P12").Value = 3

# --- Lookup table N13:P16 (Resource/Quality helper list) shifts up one row ---
$ws.Range("N13").Value = "TuTM"
$ws.Range("P13").Value = 3.5
$ws.Range("P13").VerticalAlignment = -4107

$ws.Range("N14").Value = "ThiVT"
$ws.Range("P14").Value = 4

$ws.Range("N15").Value = "HuyDV"
$ws.Range("P15").Value = 4.5

$ws.Range("N16").Value = "HuyNV"
$ws.Range("P16").Value = 5

# Row 17: H17 changes, N17 shifts up, P17 removed entirely
$ws.Range("H17").Value = 3
$ws.Range("N17").Value = "HienTT"
$ws.Range("P17").Clear()

# Row 18: G18 changes, N18 shifts up
$ws.Range("G18").Value = 3.5
$ws.Range("N18").Value = "DungDV"

# Row 19: G/H/I/J now populated (TuTM's "Insert Validate & Filter for form"), N19 shifts up
$ws.Range("G19").Value = 2
$ws.Range("G19").Font.ColorIndex = 3
$ws.Range("H19").Value = 3.5
$ws.Range("I19").Value = 3.5
$ws.Range("J19").Value = 3.5
$ws.Range("N19").Value = "HuyNV"

# Row 20 (ThiVT's task): H20 changes, N20 shifts up
$ws.Range("H20").Value = 3
$ws.Range("N20").Value = "HienTT"

# Row 21 (ThiVT's task): I21 changes, N21 shifts up
$ws.Range("I21").Value = 3.5
$ws.Range("N21").Value = "DungDV"

# Row 22 (ThiVT's task): G22 changes, N22 removed entirely
$ws.Range("G22").Value = 2
$ws.Range("G22").Font.ColorIndex = 3
$ws.Range("N22").Clear()

# Row 23: status flips to Complete, G/H/I/J now populated, K23 gets a new note
$ws.Range("F23").Value = "Complete"
$ws.Range("G23").Value = 2
$ws.Range("G23").Font.ColorIndex = 3
$ws.Range("H23").Value = 2
$ws.Range("H23").Font.ColorIndex = 3
$ws.Range("I23").Value = 3.5
$ws.Range("J23").Value = 3
$ws.Range("K23").Value = "Thiếu nhiều Bảng"

# --- Data validations: D5:D26 list range shrinks by one row (N7:N20 -> N7:N19),
#     G5:J55 list range shrinks by one row (P7:P17 -> P7:P16) ---
$ws.Range("D5:D26").Validation.Delete()
$ws.Range("D5:D26").Validation.Add(3, 1, 1, "=`$N`$7:`$N`$19")

$ws.Range("G5:J55").Validation.Delete()
$ws.Range("G5:J55").Validation.Add(3, 1, 1, "=`$P`$7:`$P`$16")
